$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.987.95"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "2.945.90"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'375.80"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").Value = "'101.24"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("E7").Value = "  -1.49%  "
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "'36.34"
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("D12").Value = "'0.0851"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "3.401.36"
$ws.Range("E13").Value = "  -1.53%  "
$ws.Range("D14").Value = "'18.13"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "'11.28"
$ws.Range("E16").Value = "  +51.77%  "
$ws.Range("D17").Value = "2.938.02"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "'0.998"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "50.945.24"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("E20").Value = "  -6.10%  "
$ws.Range("E21").Value = "  -2.94%  "
$ws.Range("D22").Value = "0.0₃0957"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").Value = "'266.63"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").Value = "'69.02"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value = "'3.18"
$ws.Range("E25").Value = "  +9.18%  "
$ws.Range("D26").Value = "'8.13"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("D27").Value = "'7.44"
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'25.69"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("E31").Value = "  -6.00%  "
$ws.Range("D32").Value = "'10.01"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").Value = "'33.36"
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("D36").Value = "'0.0442"
$ws.Range("E36").Value = "  -2.50%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "'3.13"
$ws.Range("E38").Value = "  +4.15%  "
$ws.Range("D39").Value = "'0.116"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "'16.53"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("E42").Value = "  -4.42%  "
$ws.Range("D43").Value = "'120.17"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("E44").Value = "  -2.66%  "
$ws.Range("D45").Value = "'3.42"
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D49").Value = "1.993.18"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "'0.0328"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("D51").Value = "'5.22"
$ws.Range("E51").Value = "  +1.38%  "
